# Commit: "tried to start to download resources files in case the resource
# folder does not have the right resources needed."
#
# The Bracket sheet's animal matchup cells (the "s" shared-string cells in
# columns D/E/F/G/H/I/J/K/L/M/N) are re-randomized to a different set of
# animal names. This reproduces that re-shuffle by writing the new text
# value into each affected cell on the "Bracket" worksheet.

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($s in $wb.Worksheets) {
    if ($s.Name -eq "Bracket") {
        $ws = $s
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Range('D2').Value  = 'Four-Striped Grass Mouse'
$ws.Range('N2').Value  = 'Spongilla Fly'
$ws.Range('E4').Value  = 'Four-Striped Grass Mouse'
$ws.Range('M4').Value  = 'Spongilla Fly'
$ws.Range('D6').Value  = 'Giant striped mongoose'
$ws.Range('N6').Value  = 'Puffer Fish'
$ws.Range('F8').Value  = 'Four-Striped Grass Mouse'
$ws.Range('L8').Value  = 'Spongilla Fly'
$ws.Range('D10').Value = 'Striped Possum'
$ws.Range('N10').Value = 'Trapdoor Spider'
$ws.Range('E12').Value = 'Chequered elephant shrew'
$ws.Range('M12').Value = 'Tent-making Bat'
$ws.Range('D14').Value = 'Chequered elephant shrew'
$ws.Range('N14').Value = 'Tent-making Bat'
$ws.Range('G16').Value = 'Four-Striped Grass Mouse'
$ws.Range('K16').Value = 'Spongilla Fly'
$ws.Range('D18').Value = 'Highland Streaked Tenrec'
$ws.Range('N18').Value = 'Rufous Hornero'
$ws.Range('E20').Value = 'Fire-footed Rope Squirrel'
$ws.Range('M20').Value = 'Bee'
$ws.Range('D22').Value = 'Fire-footed Rope Squirrel'
$ws.Range('N22').Value = 'Bee'
$ws.Range('F24').Value = 'Badger Bat'
$ws.Range('L24').Value = 'Dung Beetle'
$ws.Range('D26').Value = 'Numbat'
$ws.Range('N26').Value = 'New Caledonian Crow'
$ws.Range('E28').Value = 'Badger Bat'
$ws.Range('M28').Value = 'Dung Beetle'
$ws.Range('D30').Value = 'Badger Bat'
$ws.Range('N30').Value = 'Dung Beetle'
$ws.Range('H32').Value = 'Bumblebee Bat'
$ws.Range('I32').Value = 'Bumblebee Bat'
$ws.Range('J32').Value = 'Spongilla Fly'
$ws.Range('D34').Value = 'Bumblebee Bat'
$ws.Range('N34').Value = 'Lined Seahorse'
$ws.Range('C35').Value = 'Bumblebee Bat'
$ws.Range('E36').Value = 'Bumblebee Bat'
$ws.Range('M36').Value = 'Lined Seahorse'
$ws.Range('D38').Value = 'Grasshopper Mouse'
$ws.Range('N38').Value = 'Caspian Terns'
$ws.Range('F40').Value = 'Bumblebee Bat'
$ws.Range('L40').Value = 'Lined Seahorse'
$ws.Range('D42').Value = 'Silver Pika'
$ws.Range('N42').Value = 'Peacock Wrasse'
$ws.Range('E44').Value = 'Siberian Chipmunk'
$ws.Range('M44').Value = 'Darwin''s Frogs'
$ws.Range('D46').Value = 'Siberian Chipmunk'
$ws.Range('N46').Value = 'Darwin''s Frogs'
$ws.Range('G48').Value = 'Bumblebee Bat'
$ws.Range('K48').Value = 'Lined Seahorse'
$ws.Range('D50').Value = 'Silky Anteater'
$ws.Range('N50').Value = 'Spotted sandpiper'
$ws.Range('E52').Value = 'Colo Colo Opossum'
$ws.Range('M52').Value = 'Giant Water Bug'
$ws.Range('D54').Value = 'Colo Colo Opossum'
$ws.Range('N54').Value = 'Giant Water Bug'
$ws.Range('F56').Value = 'Pygmy Jerboa'
$ws.Range('L56').Value = 'Three-Spined stickleback '
$ws.Range('D58').Value = 'Thor''s Hero Shrew'
$ws.Range('N58').Value = 'Dyak Fruit Bat'
$ws.Range('E60').Value = 'Pygmy Jerboa'
$ws.Range('M60').Value = 'Three-Spined stickleback '
$ws.Range('D62').Value = 'Pygmy Jerboa'
$ws.Range('N62').Value = 'Three-Spined stickleback '
